$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1057
$ws.Range("E2").Value = 57
$ws.Range("F2").Value = 57
$ws.Range("G2").Value = 37
$ws.Range("H2").Value = 36
$ws.Range("I2").Value = 36
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1467
$ws.Range("L2").Value = 540
$ws.Range("M2").Value = 927
$ws.Range("N2").Value = 923
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 210
$ws.Range("Q2").Value = 126
$ws.Range("R2").Value = -91
$ws.Range("S2").Value = -50
$ws.Range("T2").Value = 92
$ws.Range("U2").Value = 35
$ws.Range("V2").Value = 207
$ws.Range("W2").Value = 5.39
$ws.Range("X2").Value = 3.36
$ws.Range("Y2").Value = 3.91
$ws.Range("Z2").Value = 2.44
$ws.Range("AA2").Value = 58.24
$ws.Range("AB2").Value = 339.53
$ws.Range("AC2").Value = 85
$ws.Range("AD2").Value = 25.72
$ws.Range("AE2").Value = 2197
$ws.Range("AF2").Value = 1
$ws.Range("AG2").Value = 35
$ws.Range("AH2").Value = 1.6
$ws.Range("AI2").Value = 41.1
$ws.Range("AJ2").Value = 42000000

$ws.Range("D3").Value = 1066
$ws.Range("E3").Value = 62
$ws.Range("F3").Value = 62
$ws.Range("G3").Value = 64
$ws.Range("H3").Value = 51
$ws.Range("I3").Value = 51
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1457
$ws.Range("L3").Value = 508
$ws.Range("M3").Value = 949
$ws.Range("N3").Value = 945
$ws.Range("O3").Value = 4
$ws.Range("P3").Value = 210
$ws.Range("Q3").Value = 93
$ws.Range("R3").Value = -52
$ws.Range("S3").Value = -41
$ws.Range("T3").Value = 52
$ws.Range("U3").Value = 41
$ws.Range("V3").Value = 188
$ws.Range("W3").Value = 5.81
$ws.Range("X3").Value = 4.76
$ws.Range("Y3").Value = 5.47
$ws.Range("Z3").Value = 3.47
$ws.Range("AA3").Value = 53.52
$ws.Range("AB3").Value = 350
$ws.Range("AC3").Value = 122
$ws.Range("AD3").Value = 23.04
$ws.Range("AE3").Value = 2250
$ws.Range("AF3").Value = 1.24
$ws.Range("AG3").Value = 35
$ws.Range("AH3").Value = 1.25
$ws.Range("AI3").Value = 28.8
$ws.Range("AJ3").Value = 42000000

$ws.Range("D4").Value = 1072
$ws.Range("E4").Value = 80
$ws.Range("F4").Value = 80
$ws.Range("G4").Value = 83
$ws.Range("H4").Value = 67
$ws.Range("I4").Value = 67
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1448
$ws.Range("L4").Value = 455
$ws.Range("M4").Value = 993
$ws.Range("N4").Value = 989
$ws.Range("O4").Value = 4
$ws.Range("P4").Value = 210
$ws.Range("Q4").Value = 94
$ws.Range("R4").Value = -24
$ws.Range("S4").Value = -70
$ws.Range("T4").Value = 26
$ws.Range("U4").Value = 68
$ws.Range("V4").Value = 138
$ws.Range("W4").Value = 7.48
$ws.Range("X4").Value = 6.21
$ws.Range("Y4").Value = 6.92
$ws.Range("Z4").Value = 4.58
$ws.Range("AA4").Value = 45.83
$ws.Range("AB4").Value = 371.19
$ws.Range("AC4").Value = 159
$ws.Range("AD4").Value = 18.52
$ws.Range("AE4").Value = 2355
$ws.Range("AF4").Value = 1.25
$ws.Range("AG4").Value = 35
$ws.Range("AH4").Value = 1.19
$ws.Range("AI4").Value = 21.98
$ws.Range("AJ4").Value = 42000000

$ws.Range("D5").Value = 1099
$ws.Range("E5").Value = 17
$ws.Range("F5").Value = 17
$ws.Range("G5").Value = -7
$ws.Range("H5").Value = -8
$ws.Range("I5").Value = -8
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1405
$ws.Range("L5").Value = 494
$ws.Range("M5").Value = 911
$ws.Range("N5").Value = 908
$ws.Range("O5").Value = 3
$ws.Range("P5").Value = 210
$ws.Range("Q5").Value = 75
$ws.Range("R5").Value = -58
$ws.Range("S5").Value = -17
$ws.Range("T5").Value = 50
$ws.Range("U5").Value = 25
$ws.Range("V5").Value = 210
$ws.Range("W5").Value = 1.55
$ws.Range("X5").Value = -0.77
$ws.Range("Y5").Value = -0.86
$ws.Range("Z5").Value = -0.59
$ws.Range("AA5").Value = 54.15
$ws.Range("AB5").Value = 365.97
$ws.Range("AC5").Value = -19
$ws.Range("AD5").Value = -126.63
$ws.Range("AE5").Value = 2301
$ws.Range("AF5").Value = 1.06
$ws.Range("AG5").Value = 50
$ws.Range("AH5").Value = 2.04
$ws.Range("AI5").Value = -242.86
$ws.Range("AJ5").Value = 42000000

$ws.Range("D6").Value = 1123
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = -117
$ws.Range("H6").Value = -89
$ws.Range("I6").Value = -88
$ws.Range("K6").Value = 1443
$ws.Range("L6").Value = 657
$ws.Range("M6").Value = 785
$ws.Range("N6").Value = 782
$ws.Range("P6").Value = 210
$ws.Range("Q6").Value = -2
$ws.Range("R6").Value = -105
$ws.Range("S6").Value = 121
$ws.Range("T6").Value = 94
$ws.Range("U6").Value = -96
$ws.Range("V6").Value = 364
$ws.Range("W6").Value = 0.43
$ws.Range("X6").Value = -7.91
$ws.Range("Y6").Value = -10.41
$ws.Range("Z6").Value = -6.24
$ws.Range("AA6").Value = 83.7
$ws.Range("AB6").Value = 308.57
$ws.Range("AC6").Value = -209
$ws.Range("AD6").Value = -8.67
$ws.Range("AE6").Value = 1997
$ws.Range("AF6").Value = 0.91
$ws.Range("AG6").Value = 35
$ws.Range("AH6").Value = 1.93
$ws.Range("AI6").Value = -15.59
$ws.Range("AJ6").Value = 42000000

$ws.Range("D7:AJ9").ClearContents()
